$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'245.28"
$ws.Range("G2").Formula = "'7"
$ws.Range("D3").Formula = "'22.00"
$ws.Range("G3").Formula = "'7"
$ws.Range("D4").Formula = "'5.345"
$ws.Range("G4").Formula = "'7"
$ws.Range("D5").Formula = "'0.05958"
$ws.Range("G5").Formula = "'7"
$ws.Range("D6").Formula = "'3.395"
$ws.Range("G6").Formula = "'7"
$ws.Range("D7").Formula = "'6.389"
$ws.Range("G7").Formula = "'7"
$ws.Range("D8").Formula = "'0.8102"
$ws.Range("G8").Formula = "'7"
$ws.Range("D9").Formula = "'0.9639"
$ws.Range("G9").Formula = "'7"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Formula = "'0.1425"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("G10").Formula = "'7"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Formula = "'0.07354"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("G11").Formula = "'7"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Formula = "'0.03400"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G12").Formula = "'7"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Formula = "'0.03054"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("G13").Formula = "'7"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Formula = "'0.09404"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("G14").Formula = "'7"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Formula = "'3.994"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("G15").Formula = "'7"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Formula = "'0.001596"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("G16").Formula = "'7"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Formula = "'0.04812"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("G17").Formula = "'7"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Formula = "'0.0005912"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("G18").Formula = "'7"
$ws.Range("D19").Formula = "'0.006083"
$ws.Range("G19").Formula = "'7"
$ws.Range("D20").Formula = "'0.005109"
$ws.Range("G20").Formula = "'7"
$ws.Range("D21").Formula = "'0.0009891"
$ws.Range("G21").Formula = "'7"
$ws.Range("D22").Formula = "'0.00009704"
$ws.Range("G22").Formula = "'7"
$ws.Range("D23").Formula = "'3.751"
$ws.Range("G23").Formula = "'7"
$ws.Range("G24").Formula = "'7"
$ws.Range("G25").Formula = "'7"
$ws.Range("G26").Formula = "'7"
$ws.Range("D27").Formula = "'0.0002462"
$ws.Range("G27").Formula = "'7"
$ws.Range("G28").Formula = "'7"
$ws.Range("G29").Formula = "'7"
$ws.Range("G30").Formula = "'7"
$ws.Range("G31").Formula = "'7"
$ws.Range("G32").Formula = "'7"
$ws.Range("G33").Formula = "'7"
$ws.Range("G34").Formula = "'7"
$ws.Range("G35").Formula = "'7"
$ws.Range("G36").Formula = "'7"
$ws.Range("G37").Formula = "'7"
$ws.Range("G38").Formula = "'7"
$ws.Range("G39").Formula = "'7"
$ws.Range("D40").Formula = "'0.03935"
$ws.Range("G40").Formula = "'7"
$ws.Range("D41").Formula = "'0.006373"
$ws.Range("G41").Formula = "'7"
$ws.Range("G42").Formula = "'7"
$ws.Range("D43").Formula = "'0.003001"
$ws.Range("G43").Formula = "'7"
$ws.Range("D44").Formula = "'0.005798"
$ws.Range("G44").Formula = "'7"
$ws.Range("D45").Formula = "'0.00005320"
$ws.Range("G45").Formula = "'7"
$ws.Range("G46").Formula = "'7"
$ws.Range("D47").Formula = "'0.8504"
$ws.Range("G47").Formula = "'7"
$ws.Range("D48").Formula = "'0.03423"
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("G48").Formula = "'7"
$ws.Range("G49").Formula = "'7"
$ws.Range("D50").Formula = "'0.01010"
$ws.Range("G50").Formula = "'7"
$ws.Range("G51").Formula = "'7"
